# Update the imaging-data file paths in column F (the "path" column) to
# reflect the new relative location under data/SAH/, and move the active
# selection to H10 (mirrors the author's manual edit before re-saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "./data/SAH/imaging_data/patient 1.nii.gz"
$ws.Range("F3").Value = "./data/SAH/imaging_data/patient 2.nii.gz"
$ws.Range("F4").Value = "./data/SAH/imaging_data/patient 3.nii.gz"
$ws.Range("F5").Value = "./data/SAH/imaging_data/patient 4.nii.gz"

$ws.Range("H10").Select() | Out-Null
